$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.814.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.880.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.89%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.85%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4731'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +6.46%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3978'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.21%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08060'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.23%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.027'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.91%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.36%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.894.17'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.71%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.972'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.51%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.188'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001052'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.01%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.24%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06623'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.860.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.93%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.508'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.32%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.97%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.299'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.07%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.115.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.86%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.46%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.14%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.107'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.38%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.614'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.25%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.86'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.14%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9777'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.04%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09568'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.03%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.467'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.636'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.319'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.03%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02265'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.66%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06108'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.97%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.229'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.84%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.220'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.29%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6036'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.83%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1913'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.08%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.99%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5721'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.247'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.90%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.33%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.421'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.85%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.944'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06822'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.39%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '113.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.47%  '
$ws.Range("E51").Style = "Normal"
